# Update sprint_3 oraz IN/OUT
# Uzupelnienie w module atrial_fibr wymaganych danych odnosnie sprintu

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Fill in the IN/OUT data for the ATRIAL_FIBR row (row 32)
$ws.Range("D32").Value = "interwały RR, sygnał po filtracji, QRS-onset, QRS-end, T-end, P-onset, P-end, lokalizacja R-peaks"
$ws.Range("E32").Value = "najwyższa amplituda sygnału af, częstotliwość dominująca, początek i koniec af"

# Update selection to match the author's final cursor position
$ws.Range("E33").Select()

# Turn on iterative calculation with a convergence delta of 1E-4 (calcPr/@iterateDelta)
$excel.Iteration = $true
$excel.MaxChange = 0.0001
